$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 426, shifting existing rows 426-484 down to 427-485.
$ws.Rows("426:426").Insert()

# Populate the newly inserted row 426 with the new record's data.
$ws.Range("A426").Value = 10
$ws.Range("B426").Value = "Vega Modelo de Temuco"
$ws.Range("C426").Value = "La Araucanía"
$ws.Range("D426").Value = 45131
$ws.Range("E426").Value = 9
$ws.Range("F426").Value = 100112001
$ws.Range("G426").Value = "Berenjena"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 25
$ws.Range("K426").Value = 10000
$ws.Range("L426").Value = 10000
$ws.Range("M426").Value = 10000
$ws.Range("N426").Value = "$/caja 40 unidades"
$ws.Range("O426").Value = "Región de Arica y Parinacota"
$ws.Range("P426").Value = 250
$ws.Range("Q426").Value = 40
$ws.Range("R426").Value = "Hortaliza"
